{"js": "// Apply the \"Completed Monster Details and Habits List currently adding\n// Create new Habit View\" edit described by the diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Opening line: \"Hi My name is Jacob...\" -> \"Hi, my name is Jacob...\"\n// ---------------------------------------------------------------------\nconst introResults = body.search(\"Hi My name is Jacob\", { matchCase: true });\nintroResults.load(\"items\");\nawait context.sync();\nif (introResults.items.length > 0) {\n  introResults.items[0].insertText(\"Hi, my name is Jacob\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Third paragraph (\"I've been having a lot of fun but sometimes...\")\n//    is shortened and split into two paragraphs, and the old fourth\n//    paragraph (\"This is where my application comes in...\") is replaced\n//    by the tail half of that split (with new wording).\n// ---------------------------------------------------------------------\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst funParagraph = paras.items[2];\nfunParagraph.insertText(\n  \"I\\u2019ve been having a lot of fun but sometimes it\\u2019s hard to keep track of all the numerous monsters in the game, what does each one like to do, or importantly whose hunted the biggest monster out of me and of all my friends?\",\n  \"Replace\"\n);\nawait context.sync();\n\nconst appParagraph = funParagraph.insertParagraph(\n  \"This is where my application comes in. The app lets you and your friends record all your hunt attempts on various monsters. Which means you can compare times or check on monsters\\u2019 details to prepare for their hunts.\",\n  \"After\"\n);\nawait context.sync();\n\n// The old paragraph that used to read \"This is where my application comes\n// in...combat their various moves.\" now sits right after appParagraph;\n// remove it since its content has been superseded above.\nconst parasAfterSplit = context.document.body.paragraphs;\nparasAfterSplit.load(\"items\");\nawait context.sync();\n\nconst oldAppParagraph = parasAfterSplit.items[4];\noldAppParagraph.load(\"text\");\nawait context.sync();\nif (oldAppParagraph.text.indexOf(\"This is where my application comes in\") !== -1) {\n  oldAppParagraph.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Final paragraph (\"Throughout our projects we used Agile...Now I\n//    will pass it onto Chris.\") is rewritten and expanded into five\n//    paragraphs covering Scrum, Kanban, the MVP process, and the\n//    hand-off to Ronil.\n// ---------------------------------------------------------------------\nconst parasFinal = context.document.body.paragraphs;\nparasFinal.load(\"items\");\nawait context.sync();\n\nconst agileParagraph = parasFinal.items[parasFinal.items.length - 1];\nagileParagraph.insertText(\n  \"Throughout our projects we utilised the scrum methodology. Our project was done over the course of a week, so we split our sprints into days which always started with a stand-up to check up on everyone\\u2019s progress. The daily sprints ended with retros that helped us to reflect and improve for the sprint the day after.\",\n  \"Replace\"\n);\nawait context.sync();\n\nconst followUpTexts = [\n  \"For the bigger scope we used Kanban boards on GitHub like the one you see below, which we used to manage our projects, the incremental development means we always had a demo-able application from the earliest possible stage which we would show in our stand-ups and retros. \",\n  \"The process was repeated until we had all reached our MVP.\",\n  \"I will now pass it onto Ronil to talk about the Project Process.\"\n];\n\nlet lastParagraph = agileParagraph;\nfor (const text of followUpTexts) {\n  lastParagraph = lastParagraph.insertParagraph(text, \"After\");\n  await context.sync();\n}\n", "ps1": "# Apply the \"Completed Monster Details and Habits List currently adding\n# Create new Habit View\" edit described by the diff.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParaText($para, [string]$text) {\n    $r = $para.Range\n    # Exclude the trailing paragraph mark so the assignment replaces the\n    # paragraph's whole visible content without swallowing the mark / the\n    # following paragraph.\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $text\n}\n\n# ---------------------------------------------------------------------\n# 1) Opening line: \"Hi My name is Jacob...\" -> \"Hi, my name is Jacob...\"\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Hi My name is Jacob\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Hi, my name is Jacob\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 2) Third paragraph (\"I've been having a lot of fun but sometimes...\")\n#    is shortened and split into two paragraphs, and the old fourth\n#    paragraph (\"This is where my application comes in...\") is replaced\n#    by the tail half of that split (with new wording).\n# ---------------------------------------------------------------------\n$funParaIndex = 3\n$funPara = $d.Paragraphs.Item($funParaIndex)\nSet-ParaText $funPara \"I\u2019ve been having a lot of fun but sometimes it\u2019s hard to keep track of all the numerous monsters in the game, what does each one like to do, or importantly whose hunted the biggest monster out of me and of all my friends?\"\n\n$funPara = $d.Paragraphs.Item($funParaIndex)\n$funPara.Range.InsertParagraphAfter() | Out-Null\n$appParaIndex = $funParaIndex + 1\n$appPara = $d.Paragraphs.Item($appParaIndex)\nSet-ParaText $appPara \"This is where my application comes in. The app lets you and your friends record all your hunt attempts on various monsters. Which means you can compare times or check on monsters\u2019 details to prepare for their hunts.\"\n\n# The old paragraph that used to read \"This is where my application comes\n# in...combat their various moves.\" now sits right after appPara; remove\n# it since its content has been superseded above.\n$oldAppParaIndex = $appParaIndex + 1\n$oldAppPara = $d.Paragraphs.Item($oldAppParaIndex)\nif ($oldAppPara.Range.Text -like \"*This is where my application comes in*\") {\n    $oldAppPara.Range.Delete() | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# 3) Final paragraph (\"Throughout our projects we used Agile...Now I\n#    will pass it onto Chris.\") is rewritten and expanded into five\n#    paragraphs covering Scrum, Kanban, the MVP process, and the\n#    hand-off to Ronil.\n# ---------------------------------------------------------------------\n$lastIndex = $d.Paragraphs.Count\n$agilePara = $d.Paragraphs.Item($lastIndex)\nSet-ParaText $agilePara \"Throughout our projects we utilised the scrum methodology. Our project was done over the course of a week, so we split our sprints into days which always started with a stand-up to check up on everyone\u2019s progress. The daily sprints ended with retros that helped us to reflect and improve for the sprint the day after.\"\n\n$followUpTexts = @(\n    \"For the bigger scope we used Kanban boards on GitHub like the one you see below, which we used to manage our projects, the incremental development means we always had a demo-able application from the earliest possible stage which we would show in our stand-ups and retros. \",\n    \"The process was repeated until we had all reached our MVP.\",\n    \"I will now pass it onto Ronil to talk about the Project Process.\"\n)\n\nforeach ($text in $followUpTexts) {\n    $curPara = $d.Paragraphs.Item($lastIndex)\n    $curPara.Range.InsertParagraphAfter() | Out-Null\n    $lastIndex = $lastIndex + 1\n    $newPara = $d.Paragraphs.Item($lastIndex)\n    Set-ParaText $newPara $text\n}\n"}
